# Apply updated cryptos list values (prices & 1h volume %) to match
# the target commit: refreshed quotes for the existing rows, plus a
# TRON/Chainlink row swap (rows 14-15: rank order changed upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces Excel to store numeric-looking text
    # (e.g. "0.105", "8.99") as Text rather than silently coercing it
    # to a Number, matching the source inlineStr cells.
    $cell.Value = "'" + $val
    # Re-flatten to the default style so the apostrophe-forced entry
    # doesn't leave a stray quotePrefix cell style behind.
    $cell.Style = "Normal"
}

Set-TextCell "D2" "43.219.67"
Set-TextCell "E2" "  +2.08%  "
Set-TextCell "D3" "2.372.27"
Set-TextCell "E3" "  +6.60%  "
Set-TextCell "E4" "  -0.15%  "
Set-TextCell "D5" "310.49"
Set-TextCell "E5" "  +4.28%  "
Set-TextCell "D6" "106.39"
Set-TextCell "E6" "  -5.70%  "
Set-TextCell "D7" "0.640"
Set-TextCell "E7" "  +1.85%  "
Set-TextCell "E8" "  -0.17%  "
Set-TextCell "D9" "0.635"
Set-TextCell "E9" "  +3.19%  "
Set-TextCell "D10" "43.04"
Set-TextCell "E10" "  -5.47%  "
Set-TextCell "E11" "  +1.52%  "
Set-TextCell "D12" "8.99"
Set-TextCell "E12" "  +0.75%  "
Set-TextCell "D13" "1.08"
Set-TextCell "E13" "  +13.46%  "
Set-TextCell "B14" "TRON"
Set-TextCell "C14" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D14" "0.105"
Set-TextCell "E14" "  +0.93%  "
Set-TextCell "B15" "Chainlink"
Set-TextCell "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D15" "16.50"
Set-TextCell "E15" "  +8.58%  "
Set-TextCell "D16" "2.730.11"
Set-TextCell "E16" "  +6.64%  "
Set-TextCell "D17" "2.373.46"
Set-TextCell "E17" "  +4.89%  "
Set-TextCell "D18" "43.221.89"
Set-TextCell "E18" "  +1.91%  "
Set-TextCell "E19" "  +1.93%  "
Set-TextCell "D20" "7.41"
Set-TextCell "E20" "  +0.08%  "
Set-TextCell "D21" "75.57"
Set-TextCell "D22" "3.45"
Set-TextCell "E22" "  -3.18%  "
Set-TextCell "E23" "  +8.11%  "
Set-TextCell "D24" "253.53"
Set-TextCell "E24" "  +10.04%  "
Set-TextCell "E25" "  -5.68%  "
Set-TextCell "D26" "12.07"
Set-TextCell "E26" "  +2.04%  "
Set-TextCell "E27" "  -0.07%  "
Set-TextCell "D28" "39.05"
Set-TextCell "E28" "  -0.29%  "
Set-TextCell "D29" "2.25"
Set-TextCell "E29" "  +1.36%  "
Set-TextCell "D30" "22.84"
Set-TextCell "E30" "  +7.61%  "
Set-TextCell "D31" "173.05"
Set-TextCell "E31" "  -0.64%  "
Set-TextCell "E32" "  -2.32%  "
Set-TextCell "D33" "0.0910"
Set-TextCell "E33" "  +2.04%  "
Set-TextCell "D34" "5.85"
Set-TextCell "E34" "  +0.74%  "
Set-TextCell "D35" "4.99"
Set-TextCell "E35" "  +2.02%  "
Set-TextCell "E36" "  +3.32%  "
Set-TextCell "D37" "0.0378"
Set-TextCell "E37" "  +2.30%  "
Set-TextCell "D38" "4.08"
Set-TextCell "E38" "  -5.68%  "
Set-TextCell "E39" "  +0.19%  "
Set-TextCell "E40" "  +11.30%  "
Set-TextCell "D41" "1.53"
Set-TextCell "E41" "  +15.34%  "
Set-TextCell "D42" "72.46"
Set-TextCell "E42" "  +1.39%  "
Set-TextCell "D43" "0.233"
Set-TextCell "E43" "  -2.97%  "
Set-TextCell "E44" "  -0.07%  "
Set-TextCell "D45" "12.38"
Set-TextCell "E45" "  -6.26%  "
Set-TextCell "E46" "  +3.45%  "
Set-TextCell "D47" "9.41"
Set-TextCell "E47" "  +9.59%  "
Set-TextCell "D48" "112.83"
Set-TextCell "E48" "  +7.12%  "
Set-TextCell "E49" "  -1.68%  "
Set-TextCell "D50" "0.0996"
Set-TextCell "E50" "  +1.05%  "
Set-TextCell "D51" "1.499.39"
Set-TextCell "E51" "  +4.44%  "
